$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value updates ---
$ws.Range("E1").Value = "addSequence"
$ws.Range("B2").Value = "Add Lead process workflow using  Nucleotide sequence"
$ws.Range("E2").Value = "GTCGGATGATTCAAGCTCACGGGGACGAGCAGGAGCGCTCTCGACTTTTCTAGAGCCTCAGCGTCCTAGG`nACTCACCTTTCCCTGATCCTGCACCGTCCCTCTCCTGGCCCCAGACTCTCCCTCCCACTGTTCACGAAGC`nCCAGGTGGGCCGTCGGCCGGGGAGCGGAGGGGGCGCGTGGGGTGCAGGCGGCGCCAAGGGCGCGTGCACC`nTGTGGGCGCGGGGCGCGAGGGCCCCTCCCGGCGCGAGCGGGCGCAGTTCCCCGGCGGCGCCGCTAGGGGT`nCTCTCTCGGGTGCCGAGCGGGGTGGGCCGGATCAGCTGACTCGCCTGGCTCTGAGCCCCGCCGCCGCGCT`nCGGGCTCCGTCAGTTTCCTCGGCAGCGGTAGGCGAGAGCACGCGGAGGAGCGTGCGCGGGGGCCCCGGGA`nGACGGCGGCGGTGGCGGCGCGGGCAGAGCAAGGACGCGGCGGATCCCACTCGCACAGCAGCGCACTCGGT`nGCCCCGCGCAGGGTCGCGATGCTGCCCGGTTTGGCACTGCTCCTGCTGGCCGCCTGGACGGCTCGGGCGC`nTGGAGGTGGGTGCCGCGCCTCGGAAGGCGGGGGGAGGCTGCACGGTGGGGACGCGATACCCCCCAAGACC`nTTAACCCAAGTCTTTAATGCAGAGAAGCCGGGGGTCCGTCAATGGGACCCCTCTCCTCTCCGCCCCCGCT`nTGCGGACGTCCAGCGCATCCCCGCTTTCGGCCCAGCCCTGCCCCAGGGAGTCGCGCTCCGGCCCGCTGAG`nAGGGAGCGGGCGAGGCGCTGGTCTCCCTGGTTCCGCGCCAGCCCGGGGCGAGAAGGGTAGGGGGCGACCC`nTGAGCCCAGACCCCGACTTAGTCCCTGCCTTGGAAGCGGGGGTCGGGGGAGGCGAGAGACATTCAGACAG`n"
$ws.Range("F2").Value = "Butler KM"

# --- Wrap text on E2 ---
$ws.Range("E2").WrapText = $true

# --- Row height for row 2 ---
$ws.Rows.Item(2).RowHeight = 40.5

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 51.166667
$ws.Columns.Item(4).ColumnWidth = 18.666667
$ws.Columns.Item(5).ColumnWidth = 51.0

# --- Selection ---
$ws.Range("G2").Select()
